$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# New column F: "time_taken" metadata header, styled like the other headers (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats - reuse existing header style

# Per-row timestamps recording how long each panel entry took to process
$timestamps = @(
    "2021-10-05 13:40:38.495465",
    "2021-10-05 13:40:38.495476",
    "2021-10-05 13:40:38.495481",
    "2021-10-05 13:40:38.495484",
    "2021-10-05 13:40:38.495487",
    "2021-10-05 13:40:38.495490",
    "2021-10-05 13:40:38.495493",
    "2021-10-05 13:40:38.495496",
    "2021-10-05 13:40:38.495500",
    "2021-10-05 13:40:38.495503",
    "2021-10-05 13:40:38.495506",
    "2021-10-05 13:40:38.495509",
    "2021-10-05 13:40:38.495512",
    "2021-10-05 13:40:38.495514",
    "2021-10-05 13:40:38.495517",
    "2021-10-05 13:40:38.495521"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
